# GITBOOK-240: change request with no subject merged in GitBook
#
# Reworks the "Screws, Nuts and Bolts" and "Prints" sections of the
# "Buggy components list" table: new fastener rows are inserted, a
# couple of quantities change, and the bottom of the "Prints" section
# grows by three new printed parts. The backing Excel Table
# (CustomerList) grows from B4:F59 to B4:F62 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# "Screws, Nuts and Bolts" section (rows 31-45).
# Row 30 stays blank (divider). Row 31 is a new row ("M2-12mm Bolt")
# that reuses the same look as the blank divider rows (style already
# present on B30:F31), so only values need to be written.
# ---------------------------------------------------------------
$fasteners = @(
    @{ Row = 31; Name = "M2-12mm Bolt";        Qty = 1  },
    @{ Row = 32; Name = "M3-12mm Bolt";        Qty = 44 },
    @{ Row = 33; Name = "M3-5mm Bolt";         Qty = 4  },
    @{ Row = 34; Name = "M3-25mm Bolt";        Qty = 2  },
    @{ Row = 35; Name = "M3-20mm Bolt";        Qty = 6  },
    @{ Row = 36; Name = "M3-40mm Bolt";        Qty = 1  },
    @{ Row = 37; Name = "M3.5-10mm Bolt";      Qty = 6  },
    @{ Row = 38; Name = "M3.5-35mm Bolt";      Qty = 8  },
    @{ Row = 39; Name = "M5-20mm Bolt";        Qty = 2  },
    @{ Row = 40; Name = "M2 nut";              Qty = 1  },
    @{ Row = 41; Name = "M3 nut";              Qty = 16 },
    @{ Row = 42; Name = "M3.5 nut";            Qty = 8  },
    @{ Row = 43; Name = "M5 nut";              Qty = 2  },
    @{ Row = 44; Name = "M3-10mm nylon Screw"; Qty = 4  },
    @{ Row = 45; Name = "M3 nylon nut";        Qty = 4  }
)

foreach ($item in $fasteners) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Name
    $ws.Range("C$r").Value = $item.Qty
    $ws.Range("E$r").Value = "any store"
}

# ---------------------------------------------------------------
# "Prints" section header stays put (now row 46); the first item
# (BottomPlate, row 47) keeps its distinctive style, the rest
# (rows 48-62) are plain rows that already carry the right style -
# three new parts (CameraFrame, Pixy2Box1, Pixy2Box2) are appended
# at the end, reusing the previously-blank striped rows 60-62.
# ---------------------------------------------------------------
$prints = @(
    @{ Row = 47; Name = "BottomPlate";              Qty = 1 },
    @{ Row = 48; Name = "UpperPlate";                Qty = 1 },
    @{ Row = 49; Name = "WheelPlate";                Qty = 1 },
    @{ Row = 50; Name = "FrontWheelAxle";            Qty = 2 },
    @{ Row = 51; Name = "ServoArm";                  Qty = 2 },
    @{ Row = 52; Name = "FrontWheelBearingHolder";   Qty = 2 },
    @{ Row = 53; Name = "Fuzeta";                    Qty = 2 },
    @{ Row = 54; Name = "LowerClamp";                Qty = 4 },
    @{ Row = 55; Name = "UpperClamp";                Qty = 4 },
    @{ Row = 56; Name = "DisplayBox";                Qty = 1 },
    @{ Row = 57; Name = "CameraHold";                Qty = 1 },
    @{ Row = 58; Name = "SlidingComponent";          Qty = 1 },
    @{ Row = 59; Name = "RotatingComponent";         Qty = 1 },
    @{ Row = 60; Name = "CameraFrame";                Qty = 1 },
    @{ Row = 61; Name = "Pixy2Box1";                  Qty = 1 },
    @{ Row = 62; Name = "Pixy2Box2";                  Qty = 1 }
)

foreach ($item in $prints) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Name
    $ws.Range("C$r").Value = $item.Qty
}

# Rows 60-62 used to be blank striped filler rows (style only, no
# values) below the table; now that they hold real data they must
# pick up the same formatting as the other plain "Prints" rows (e.g.
# row 59) instead of the old striped-filler look.
$ws.Range("B59:F59").Copy() | Out-Null
$ws.Range("B60:F62").PasteSpecial(-4122) | Out-Null

# The filler/striping below the table shifts down by three rows (it
# used to start at row 60, now it starts at row 63) but keeps
# alternating between the same two fills.
$ws.Range("B60:F61").Copy() | Out-Null
$ws.Range("B63:F70").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Grow the backing table (CustomerList) to cover the three new rows.
# ---------------------------------------------------------------
$lo = $ws.ListObjects.Item("CustomerList")
$lo.Resize($ws.Range("B4:F62")) | Out-Null

# ---------------------------------------------------------------
# Restore the view: scrolled down a bit further and a different
# active cell than when the edit started.
# ---------------------------------------------------------------
$ws.Range("E35").Select() | Out-Null
